$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row-level changes ---------------------------------------------------
# Two data rows were removed entirely: "RM 232" (row 26) and "SC 92"
# (originally row 28, which becomes row 27 once row 26 is deleted).
# Deleting shifts every following row up, matching the new dimension A1:F33.
$ws.Rows("26").Delete()
$ws.Rows("27").Delete()

# --- Cell-level value changes (rows 2-25, unaffected by the deletion) ----
$ws.Range("E2").Value = ""

$ws.Range("E5").Value = -5

$ws.Range("C6").Value = 15.1
$ws.Range("E6").Value = -5.7

$ws.Range("C8").Value = ""

$ws.Range("E10").Value = ""

$ws.Range("C12").Value = 12.5

$ws.Range("E13").Value = ""

$ws.Range("C14").Value = ""

$ws.Range("C17").Value = 11.2

$ws.Range("C18").Value = 11.5

$ws.Range("C19").Value = ""

$ws.Range("C20").Value = ""

$ws.Range("C23").Value = 12.2

$ws.Range("E24").Value = -8.1

# --- Cell-level value changes (rows 26-33, after the two rows above ------
# --- were deleted and everything shifted up) ------------------------------
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = ""

$ws.Range("E28").Value = ""

$ws.Range("B29").Value = ""

$ws.Range("E30").Value = -5.7

$ws.Range("B32").Value = ""
